# Shift a set of shapes on slide 1 horizontally (only the x offset of
# <a:off>; y and the <a:ext> size stay untouched) to match the updated
# "background" layout.
#
# The PowerPoint object model here stores Shape.Left/.Top as single
# precision (float32) point values, and converts back to EMU (the unit
# used in the underlying OOXML <a:off x="..."/>) by truncating
# floor(pointsAsFloat32 * 12700). A naive `target_emu / 12700` assignment
# can therefore round-trip to `target_emu - 1` once it is re-quantized to
# float32. To hit the exact EMU value requested by the diff we search,
# in small float increments, for a point value whose float32 rounding
# reproduces the desired EMU exactly.

function Get-EmuFromPoints($ptValue) {
    $f = [float]$ptValue
    $emu = [double]$f * 12700
    return [math]::Floor($emu)
}

function Find-PointsForEmu($targetEmu) {
    $basePt = $targetEmu / 12700.0
    if ((Get-EmuFromPoints $basePt) -eq $targetEmu) { return $basePt }
    for ($i = 1; $i -le 20000; $i++) {
        $delta = $i * 0.0000001
        $tryUp = $basePt + $delta
        if ((Get-EmuFromPoints $tryUp) -eq $targetEmu) { return $tryUp }
        $tryDown = $basePt - $delta
        if ((Get-EmuFromPoints $tryDown) -eq $targetEmu) { return $tryDown }
    }
    throw "Find-PointsForEmu: no solution found for target EMU $targetEmu"
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Map of shape name -> new x offset, expressed in EMU (matches the target
# OOXML <a:off x="..."/> values). Top (y) is unchanged for every shape.
$newLeftEmu = @{
    "Rounded Rectangle 11"        = 5576018
    "Rectangle 12"                = 6379610
    "Rectangle 13"                = 6532010
    "Rectangle 14"                = 6684410
    "TextBox 15"                  = 6684410
    "Rectangle 16"                = 11323029
    "Rectangle 17"                = 11475429
    "Rectangle 18"                = 11627829
    "TextBox 19"                  = 11627829
    "TextBox 21"                  = 5422752
    "Rectangle 25"                = 1524662
    "Rectangle 26"                = 1677062
    "Rectangle 27"                = 1829462
    "TextBox 28"                  = 1829462
    "Oval 38"                     = 10849925
    "TextBox 39"                  = 11195325
    "Straight Arrow Connector 40" = 12593196
    "Rounded Rectangle 48"        = 11634506
    "TextBox 49"                  = 11543161
    "Right Arrow 33"              = 4949100
    "Right Arrow 54"              = 9836901
    "Right Arrow 57"              = 9836901
}

foreach ($name in $newLeftEmu.Keys) {
    $shape = $s.Shapes.Item($name)
    $targetEmu = $newLeftEmu[$name]
    $shape.Left = Find-PointsForEmu $targetEmu
}
